# Update test case wording across Sheet2, Sheet3 and Sheet4 as part of the
# "Last update before transfer" commit.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet2 - "Verify sorting functionality works correctly with applied filters"
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Sheet2")

$ws2.Range("I2").Value = "User is on the transaction list page with filters applied.`nTransaction list remains sorted and filtered as per user actions."
$ws2.Range("I3").Value = "User is on the transaction list page with filters applied.`nTransaction list remains sorted and filtered as per user actions."

# ---------------------------------------------------------------------------
# Sheet3 - "Verify the help guide is updated with filtering instructions"
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Sheet3")

$ws3.Range("B2").Value = "Verify help guide is updated with filtering instructions"
$ws3.Range("D2").Value = "Open the help guide from the transaction list page."
$ws3.Range("F2").Value = "Help guide is displayed."
$ws3.Range("I2").Value = "User has access to the help guide.`nUser is informed about filtering functionality through the help guide."

$ws3.Range("B3").Value = "Verify help guide is updated with filtering instructions"
$ws3.Range("I3").Value = "User has access to the help guide.`nUser is informed about filtering functionality through the help guide."

# ---------------------------------------------------------------------------
# Sheet4 - "Verify that the transaction graph does not refresh based on ... filters"
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("Sheet4")

$ws4.Range("B2").Value = "Verify that the transaction graph does not refresh based on the chosen filters"
$ws4.Range("I2").Value = "User is on the transaction list page with filters applied.`nTransaction graph remains static regardless of list filters."

$ws4.Range("B3").Value = "Verify that the transaction graph does not refresh based on the chosen filters"
$ws4.Range("F3").Value = "The transaction graph remains unchanged and does not refresh based on the applied filters."
$ws4.Range("I3").Value = "User is on the transaction list page with filters applied.`nTransaction graph remains static regardless of list filters."
